$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2076.4443
$ws.Range("J17").Value = 2542.6667
$ws.Range("L17").Value = 7628.000100000001
$ws.Range("N17").Value = -7964.000100000001

# Row 113
$ws.Range("H113").Value = 3881.6428
$ws.Range("I113").Value = 3874.5
$ws.Range("J113").Value = 3924.5
$ws.Range("K113").Value = 3874.5
$ws.Range("L113").Value = 3924.5
$ws.Range("M113").Value = -620.5
$ws.Range("N113").Value = -10432.5

# Row 123
$ws.Range("H123").Value = 69029.60000000001
$ws.Range("J123").Value = 69029.60000000001
$ws.Range("L123").Value = 69029.60000000001
$ws.Range("N123").Value = -78829.60000000001

# Row 132
$ws.Range("H132").Value = 2224190.8
$ws.Range("I132").Value = 1904.561
$ws.Range("J132").Value = 25002624
$ws.Range("K132").Value = 5713.683
$ws.Range("L132").Value = 75007872
$ws.Range("M132").Value = -3183.683
$ws.Range("N132").Value = -75012932

# Row 137
$ws.Range("H137").Value = 750905
$ws.Range("I137").Value = 1074.1875
$ws.Range("K137").Value = 3222.5625
$ws.Range("M137").Value = -672.5625

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1216.5416
$ws.Range("I2").Value = 467.16666
$ws.Range("J2").Value = 3464.6667
$ws.Range("K2").Value = 467.16666
$ws.Range("L2").Value = 3464.6667
$ws.Range("M2").Value = -354.16666
$ws.Range("N2").Value = -3690.6667

# Row 32
$ws.Range("H32").Value = 152269.86
$ws.Range("I32").Value = 152269.86
$ws.Range("K32").Value = 152269.86
$ws.Range("M32").Value = -151982.86

# Row 61
$ws.Range("H61").Value = 300263.84
$ws.Range("I61").Value = 2709.3171
$ws.Range("K61").Value = 2709.3171
$ws.Range("M61").Value = -2497.3171

# Row 116
$ws.Range("H116").Value = 1216.5416
$ws.Range("I116").Value = 467.16666
$ws.Range("J116").Value = 3464.6667
$ws.Range("K116").Value = 467.16666
$ws.Range("L116").Value = 3464.6667
$ws.Range("M116").Value = 1826.83334
$ws.Range("N116").Value = -8052.6667

# Row 125
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

# Row 132
$ws.Range("H132").Value = 1462.5161
$ws.Range("I132").Value = 1392.5
$ws.Range("J132").Value = 1826.6
$ws.Range("K132").Value = 4177.5
$ws.Range("L132").Value = 5479.799999999999
$ws.Range("M132").Value = -1647.5
$ws.Range("N132").Value = -10539.8

# Row 136
$ws.Range("H136").Value = 300263.84
$ws.Range("I136").Value = 2709.3171
$ws.Range("K136").Value = 8127.951300000001
$ws.Range("M136").Value = -5577.951300000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1216.5416
$ws.Range("I3").Value = 467.16666
$ws.Range("J3").Value = 3464.6667
$ws.Range("K3").Value = 467.16666
$ws.Range("L3").Value = 3464.6667
$ws.Range("M3").Value = -353.16666
$ws.Range("N3").Value = -3692.6667

# Row 19
$ws.Range("H19").Value = 12899
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 40
$ws.Range("H40").Value = 75000
$ws.Range("J40").Value = 75000
$ws.Range("L40").Value = 75000
$ws.Range("N40").Value = -75530

# Row 107
$ws.Range("H107").Value = 4596.278
$ws.Range("I107").Value = 2737.3333
$ws.Range("K107").Value = 2737.3333
$ws.Range("M107").Value = -817.3332999999998

# Row 134
$ws.Range("H134").Value = 570191.4399999999
$ws.Range("I134").Value = 2056.3489
$ws.Range("K134").Value = 6169.0467
$ws.Range("M134").Value = -3634.0467

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3306.6785
$ws.Range("I31").Value = 1218.3334
$ws.Range("J31").Value = 4872.9375
$ws.Range("K31").Value = 1218.3334
$ws.Range("L31").Value = 4872.9375
$ws.Range("M31").Value = -923.3334
$ws.Range("N31").Value = -5462.9375

# Row 34
$ws.Range("H34").Value = 3306.6785
$ws.Range("I34").Value = 1218.3334
$ws.Range("J34").Value = 4872.9375
$ws.Range("K34").Value = 1218.3334
$ws.Range("L34").Value = 4872.9375
$ws.Range("M34").Value = -1016.3334
$ws.Range("N34").Value = -5276.9375

# Row 50
$ws.Range("H50").Value = 47499
$ws.Range("J50").Value = 49999
$ws.Range("L50").Value = 49999
$ws.Range("N50").Value = -51249

# Row 51
$ws.Range("H51").Value = 20015

# Row 59
$ws.Range("H59").Value = 67814.57000000001
$ws.Range("J59").Value = 67814.57000000001
$ws.Range("L59").Value = 67814.57000000001
$ws.Range("N59").Value = -70104.57000000001

# Row 61
$ws.Range("H61").Value = 20015

# Row 68
$ws.Range("H68").Value = 75000
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76498

# Row 71
$ws.Range("H71").Value = 75000
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -232488

# Row 107
$ws.Range("H107").Value = 413.53333
$ws.Range("I107").Value = 292.5
$ws.Range("J107").Value = 551.8570999999999
$ws.Range("K107").Value = 292.5
$ws.Range("L107").Value = 551.8570999999999
$ws.Range("M107").Value = 1627.5
$ws.Range("N107").Value = -4391.8571

# Row 132
$ws.Range("H132").Value = 1815.8823
$ws.Range("I132").Value = 1211.1724
$ws.Range("J132").Value = 5323.2
$ws.Range("K132").Value = 3633.5172
$ws.Range("L132").Value = 15969.6
$ws.Range("M132").Value = -1103.5172
$ws.Range("N132").Value = -21029.6

$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 351.8421
$ws.Range("J41").Value = 1089.8334
$ws.Range("L41").Value = 3269.5002
$ws.Range("N41").Value = -3945.5002

# Row 70
$ws.Range("H70").Value = 6286.5713
$ws.Range("I70").Value = 2012
$ws.Range("K70").Value = 6036
$ws.Range("M70").Value = -5721

# Row 73
$ws.Range("H73").Value = 6286.5713
$ws.Range("I73").Value = 2012
$ws.Range("K73").Value = 6036
$ws.Range("M73").Value = -4944

# Row 113
$ws.Range("H113").Value = 1166.6154
$ws.Range("I113").Value = 1562
$ws.Range("J113").Value = 827.7143
$ws.Range("K113").Value = 4686
$ws.Range("L113").Value = 2483.1429
$ws.Range("M113").Value = -2516
$ws.Range("N113").Value = -6823.1429

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2162.6924
$ws.Range("I22").Value = 1132.8334
$ws.Range("J22").Value = 2471.65
$ws.Range("K22").Value = 1132.8334
$ws.Range("L22").Value = 2471.65
$ws.Range("M22").Value = -837.8334
$ws.Range("N22").Value = -3061.65

# Row 27
$ws.Range("H27").Value = 2162.6924
$ws.Range("I27").Value = 1132.8334
$ws.Range("J27").Value = 2471.65
$ws.Range("K27").Value = 1132.8334
$ws.Range("L27").Value = 2471.65
$ws.Range("M27").Value = -1025.8334
$ws.Range("N27").Value = -2685.65

# Row 132
$ws.Range("H132").Value = 3543.0293
$ws.Range("I132").Value = 2584.516
$ws.Range("K132").Value = 7753.548000000001
$ws.Range("M132").Value = -5223.548000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 17999.5
$ws.Range("I51").Value = 12000
$ws.Range("K51").Value = 12000
$ws.Range("M51").Value = -11490

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 132
$ws.Range("H132").Value = 27780692
$ws.Range("I132").Value = 30305844
$ws.Range("J132").Value = 4005
$ws.Range("K132").Value = 90917532
$ws.Range("L132").Value = 12015
$ws.Range("M132").Value = -90915002
$ws.Range("N132").Value = -17075
